$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text change (shared across Overview + zh-cn + de-de sheets):
#    "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# Helper color/underline to match the workbook's existing "HyperLink" look
# (font color FF6495ED == RGB(100,149,237), underlined)
# ---------------------------------------------------------------------------
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback File" (G)
#    for rows 2 and 3, and set "Latest Handback DateTime" (H)
# ---------------------------------------------------------------------------
$mdName = "a4dfe2e2-f6d7-4974-94d7-b6db15e1f8ee.md"
$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/c875e24a96dbbe3ed00c76d2fcde3fca83d7b064/e2e/a4dfe2e2-f6d7-4974-94d7-b6db15e1f8ee.md"

$zhXlfName = "a4dfe2e2-f6d7-4974-94d7-b6db15e1f8ee.ce64af4acec5733debc1b80b868ba3fc2ea2a474.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75dfbe27f0cc8f50f43a8c9cae6dbe5953b81520/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/a4dfe2e2-f6d7-4974-94d7-b6db15e1f8ee.ce64af4acec5733debc1b80b868ba3fc2ea2a474.zh-cn.xlf"

foreach ($r in 2, 3) {
    $fCell = $wsZh.Range("F$r")
    $wsZh.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdName) | Out-Null
    $fCell.Font.Underline = $true
    $fCell.Font.Color = $hyperlinkColor

    $gCell = $wsZh.Range("G$r")
    $wsZh.Hyperlinks.Add($gCell, $zhXlfUrl, "", "", $zhXlfName) | Out-Null
    $gCell.Font.Underline = $true
    $gCell.Font.Color = $hyperlinkColor
}

$wsZh.Range("H2").Value = "2016-03-23 07:13:36"
$wsZh.Range("H3").Value = "2016-03-23 07:13:36"

# ---------------------------------------------------------------------------
# 3) de-de sheet: populate "Latest Target File" (F) / "Latest Handback File" (G)
#    for rows 2 and 3, and set "Latest Handback DateTime" (H)
# ---------------------------------------------------------------------------
$deXlfName = "a4dfe2e2-f6d7-4974-94d7-b6db15e1f8ee.ce64af4acec5733debc1b80b868ba3fc2ea2a474.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/882b55fcbc283e02c68734904f7aed0d64cec11c/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/a4dfe2e2-f6d7-4974-94d7-b6db15e1f8ee.ce64af4acec5733debc1b80b868ba3fc2ea2a474.de-de.xlf"

foreach ($r in 2, 3) {
    $fCell = $wsDe.Range("F$r")
    $wsDe.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdName) | Out-Null
    $fCell.Font.Underline = $true
    $fCell.Font.Color = $hyperlinkColor

    $gCell = $wsDe.Range("G$r")
    $wsDe.Hyperlinks.Add($gCell, $deXlfUrl, "", "", $deXlfName) | Out-Null
    $gCell.Font.Underline = $true
    $gCell.Font.Color = $hyperlinkColor
}

$wsDe.Range("H2").Value = "2016-03-23 07:13:51"
$wsDe.Range("H3").Value = "2016-03-23 07:13:51"
